$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting (bold font, borders, alignment) from H1 header cell
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data values for I and J columns (rows 2-47)
$iValues = @{
    2 = 4; 3 = 7; 4 = 6; 5 = 7; 6 = 7; 7 = 6; 8 = 5; 9 = 8; 10 = 8;
    11 = 8; 12 = 8; 13 = 7; 14 = 7; 15 = 8; 16 = 7; 17 = 8; 18 = 8; 19 = 11; 20 = 6;
    21 = 7; 22 = 8; 23 = 7; 24 = 8; 25 = 7; 26 = 8; 27 = 9; 28 = 7; 29 = 8; 30 = 7;
    31 = 10; 32 = 6; 33 = 8; 34 = 6; 35 = 7; 36 = 7; 37 = 8; 38 = 7; 39 = 8; 40 = 8;
    41 = 6; 42 = 4; 43 = 7; 44 = 7; 45 = 6; 46 = 3; 47 = 6
}
$jValues = @{
    2 = 4; 3 = 7; 4 = 6; 5 = 7; 6 = 7; 7 = 6; 8 = 5; 9 = 8; 10 = 8;
    11 = 8; 12 = 8; 13 = 7; 14 = 7; 15 = 8; 16 = 7; 17 = 8; 18 = 8; 19 = 11; 20 = 6;
    21 = 7; 22 = 8; 23 = 7; 24 = 8; 25 = 8; 26 = 8; 27 = 9; 28 = 8; 29 = 8; 30 = 7;
    31 = 10; 32 = 7; 33 = 8; 34 = 7; 35 = 7; 36 = 7; 37 = 8; 38 = 7; 39 = 8; 40 = 8;
    41 = 7; 42 = 4; 43 = 7; 44 = 7; 45 = 6; 46 = 3; 47 = 6
}

foreach ($r in 2..47) {
    $ws.Cells.Item($r, 9).Value = $iValues[$r]
    $ws.Cells.Item($r, 10).Value = $jValues[$r]
}
